# Insert a new record row at row 213 of the "Ají" price sheet (weekly update:
# a new observation is added at the top of this date-ordered block, pushing
# the existing rows 213-272 down to 214-273).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 213:272 down to 214:273, leaving a blank row 213 (format is
# inherited from row 212 immediately above, which already carries the date
# number-format on column D).
$ws.Rows("213:213").Insert()

# Populate the newly inserted row 213 with the new record.
$ws.Cells.Item(213, 1).Value2  = 4
$ws.Cells.Item(213, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(213, 3).Value2  = "Los Lagos"
$ws.Cells.Item(213, 4).Value2  = 44736
$ws.Cells.Item(213, 5).Value2  = 10
$ws.Cells.Item(213, 6).Value2  = 100112021
$ws.Cells.Item(213, 7).Value2  = "Ají"
$ws.Cells.Item(213, 8).Value2  = "Inferno"
$ws.Cells.Item(213, 9).Value2  = "Primera"
$ws.Cells.Item(213, 10).Value2 = 160
$ws.Cells.Item(213, 11).Value2 = 22000
$ws.Cells.Item(213, 12).Value2 = 28000
$ws.Cells.Item(213, 13).Value2 = 25000
$ws.Cells.Item(213, 14).Value2 = "$/caja 12 kilos"
$ws.Cells.Item(213, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(213, 16).Value2 = 2083
$ws.Cells.Item(213, 17).Value2 = 12
$ws.Cells.Item(213, 18).Value2 = "Hortaliza"
